$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 7).Value = 6.488422666666668
$ws.Cells.Item(2, 8).Value = 19.465268
$ws.Cells.Item(2, 9).Value = 0.3444204430827928
$ws.Cells.Item(2, 10).Value = 0.361825925629615
$ws.Cells.Item(2, 13).Value = 1.684496
$ws.Cells.Item(2, 14).Value = 5.053488
$ws.Cells.Item(2, 15).Value = 0.6423607101334534
$ws.Cells.Item(2, 16).Value = 0.7282461611889918
$ws.Cells.Item(2, 17).Value = 10.92972202830934
$ws.Cells.Item(2, 18).Value = 98.36749825478401
$ws.Cells.Item(2, 19).Value = 0.2212421604031415
$ws.Cells.Item(2, 20).Value = 0.2634983413584207

# Row 3
$ws.Cells.Item(3, 7).Value = 6.488422666666668
$ws.Cells.Item(3, 8).Value = 19.465268
$ws.Cells.Item(3, 9).Value = 0.3444204430827928
$ws.Cells.Item(3, 10).Value = 0.361825925629615
$ws.Cells.Item(3, 13).Value = 0.9277985
$ws.Cells.Item(3, 14).Value = 1.855597
$ws.Cells.Item(3, 15).Value = 0.3538039290807178
$ws.Cells.Item(3, 16).Value = 0.2674056793968462
$ws.Cells.Item(3, 17).Value = 6.019948817499334
$ws.Cells.Item(3, 18).Value = 36.119692904996
$ws.Cells.Item(3, 19).Value = 0.1218573060184139
$ws.Cells.Item(3, 20).Value = 0.09675430746637995

# Row 4
$ws.Cells.Item(4, 7).Value = 6.488422666666668
$ws.Cells.Item(4, 8).Value = 19.465268
$ws.Cells.Item(4, 9).Value = 0.3444204430827928
$ws.Cells.Item(4, 10).Value = 0.361825925629615
$ws.Cells.Item(4, 13).Value = 0.01005766666666667
$ws.Cells.Item(4, 14).Value = 0.030173
$ws.Cells.Item(4, 15).Value = 0.003835360785828855
$ws.Cells.Item(4, 16).Value = 0.004348159414162149
$ws.Cells.Item(4, 17).Value = 0.06525839237377778
$ws.Cells.Item(4, 18).Value = 0.587325531364
$ws.Cells.Item(4, 19).Value = 0.001320976661237543
$ws.Cells.Item(4, 20).Value = 0.001573276804814344

# Row 5
$ws.Cells.Item(5, 9).Value = 0.4517209651039303
$ws.Cells.Item(5, 10).Value = 0.4745489404232121
$ws.Cells.Item(5, 13).Value = 1.684496
$ws.Cells.Item(5, 14).Value = 5.053488
$ws.Cells.Item(5, 15).Value = 0.6423607101334534
$ws.Cells.Item(5, 16).Value = 0.7282461611889918
$ws.Cells.Item(5, 17).Value = 14.33476055821333
$ws.Cells.Item(5, 18).Value = 129.01284502392
$ws.Cells.Item(5, 19).Value = 0.2901677999263296
$ws.Cells.Item(5, 20).Value = 0.3455884441595078

# Row 6
$ws.Cells.Item(6, 9).Value = 0.4517209651039303
$ws.Cells.Item(6, 10).Value = 0.4745489404232121
$ws.Cells.Item(6, 13).Value = 0.9277985
$ws.Cells.Item(6, 14).Value = 1.855597
$ws.Cells.Item(6, 15).Value = 0.3538039290807178
$ws.Cells.Item(6, 16).Value = 0.2674056793968462
$ws.Cells.Item(6, 17).Value = 7.895399777600834
$ws.Cells.Item(6, 18).Value = 47.372398665605
$ws.Cells.Item(6, 19).Value = 0.1598206523019044
$ws.Cells.Item(6, 20).Value = 0.1268970818209225

# Row 7
$ws.Cells.Item(7, 9).Value = 0.4517209651039303
$ws.Cells.Item(7, 10).Value = 0.4745489404232121
$ws.Cells.Item(7, 13).Value = 0.01005766666666667
$ws.Cells.Item(7, 14).Value = 0.030173
$ws.Cells.Item(7, 15).Value = 0.003835360785828855
$ws.Cells.Item(7, 16).Value = 0.004348159414162149
$ws.Cells.Item(7, 17).Value = 0.08558894971611111
$ws.Cells.Item(7, 18).Value = 0.770300547445
$ws.Cells.Item(7, 19).Value = 0.001732512875696379
$ws.Cells.Item(7, 20).Value = 0.002063414442781862

# Row 8
$ws.Cells.Item(8, 5).Value = 2
$ws.Cells.Item(8, 6).Value = 0.6666666666666666
$ws.Cells.Item(8, 7).Value = 0.331724
$ws.Cells.Item(8, 8).Value = 0.9951719999999999
$ws.Cells.Item(8, 9).Value = 0.01760867516355742
$ws.Cells.Item(8, 10).Value = 0.0184985395557192
$ws.Cells.Item(8, 13).Value = 1.684496
$ws.Cells.Item(8, 14).Value = 5.053488
$ws.Cells.Item(8, 15).Value = 0.6423607101334534
$ws.Cells.Item(8, 16).Value = 0.7282461611889918
$ws.Cells.Item(8, 17).Value = 0.5587877511039999
$ws.Cells.Item(8, 18).Value = 5.029089759935999
$ws.Cells.Item(8, 19).Value = 0.01131112108257205
$ws.Cells.Item(8, 20).Value = 0.01347149041905522

# Row 9
$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 0.6666666666666666
$ws.Cells.Item(9, 7).Value = 0.331724
$ws.Cells.Item(9, 8).Value = 0.9951719999999999
$ws.Cells.Item(9, 9).Value = 0.01760867516355742
$ws.Cells.Item(9, 10).Value = 0.0184985395557192
$ws.Cells.Item(9, 13).Value = 0.9277985
$ws.Cells.Item(9, 14).Value = 1.855597
$ws.Cells.Item(9, 15).Value = 0.3538039290807178
$ws.Cells.Item(9, 16).Value = 0.2674056793968462
$ws.Cells.Item(9, 17).Value = 0.3077730296139999
$ws.Cells.Item(9, 18).Value = 1.846638177684
$ws.Cells.Item(9, 19).Value = 0.006230018458772666
$ws.Cells.Item(9, 20).Value = 0.004946614537746526

# Row 10
$ws.Cells.Item(10, 5).Value = 2
$ws.Cells.Item(10, 6).Value = 0.6666666666666666
$ws.Cells.Item(10, 7).Value = 0.331724
$ws.Cells.Item(10, 8).Value = 0.9951719999999999
$ws.Cells.Item(10, 9).Value = 0.01760867516355742
$ws.Cells.Item(10, 10).Value = 0.0184985395557192
$ws.Cells.Item(10, 13).Value = 0.01005766666666667
$ws.Cells.Item(10, 14).Value = 0.030173
$ws.Cells.Item(10, 15).Value = 0.003835360785828855
$ws.Cells.Item(10, 16).Value = 0.004348159414162149
$ws.Cells.Item(10, 17).Value = 0.003336369417333333
$ws.Cells.Item(10, 18).Value = 0.030027324756
$ws.Cells.Item(10, 19).Value = 0.00006753562221270662
$ws.Cells.Item(10, 20).Value = 0.00008043459891745133

# Row 11
$ws.Cells.Item(11, 7).Value = 2.7186785
$ws.Cells.Item(11, 8).Value = 5.437357
$ws.Cells.Item(11, 9).Value = 0.1443137264130649
$ws.Cells.Item(11, 10).Value = 0.1010711349827635
$ws.Cells.Item(11, 13).Value = 1.684496
$ws.Cells.Item(11, 14).Value = 5.053488
$ws.Cells.Item(11, 15).Value = 0.6423607101334534
$ws.Cells.Item(11, 16).Value = 0.7282461611889918
$ws.Cells.Item(11, 17).Value = 4.579603058536001
$ws.Cells.Item(11, 18).Value = 27.477618351216
$ws.Cells.Item(11, 19).Value = 0.09270146778070129
$ws.Cells.Item(11, 20).Value = 0.07360466605821191

# Row 12
$ws.Cells.Item(12, 7).Value = 2.7186785
$ws.Cells.Item(12, 8).Value = 5.437357
$ws.Cells.Item(12, 9).Value = 0.1443137264130649
$ws.Cells.Item(12, 10).Value = 0.1010711349827635
$ws.Cells.Item(12, 13).Value = 0.9277985
$ws.Cells.Item(12, 14).Value = 1.855597
$ws.Cells.Item(12, 15).Value = 0.3538039290807178
$ws.Cells.Item(12, 16).Value = 0.2674056793968462
$ws.Cells.Item(12, 17).Value = 2.52238583428225
$ws.Cells.Item(12, 18).Value = 10.089543337129
$ws.Cells.Item(12, 19).Value = 0.05105876342522214
$ws.Cells.Item(12, 20).Value = 0.02702699551747622

# Row 13
$ws.Cells.Item(13, 7).Value = 2.7186785
$ws.Cells.Item(13, 8).Value = 5.437357
$ws.Cells.Item(13, 9).Value = 0.1443137264130649
$ws.Cells.Item(13, 10).Value = 0.1010711349827635
$ws.Cells.Item(13, 13).Value = 0.01005766666666667
$ws.Cells.Item(13, 14).Value = 0.030173
$ws.Cells.Item(13, 15).Value = 0.003835360785828855
$ws.Cells.Item(13, 16).Value = 0.004348159414162149
$ws.Cells.Item(13, 17).Value = 0.02734356212683333
$ws.Cells.Item(13, 18).Value = 0.164061372761
$ws.Cells.Item(13, 19).Value = 0.000553495207141503
$ws.Cells.Item(13, 20).Value = 0.0004394734070753562

# Row 14
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 0.790022
$ws.Cells.Item(14, 8).Value = 2.370066
$ws.Cells.Item(14, 9).Value = 0.04193619023665445
$ws.Cells.Item(14, 10).Value = 0.04405545940869034
$ws.Cells.Item(14, 13).Value = 1.684496
$ws.Cells.Item(14, 14).Value = 5.053488
$ws.Cells.Item(14, 15).Value = 0.6423607101334534
$ws.Cells.Item(14, 16).Value = 0.7282461611889918
$ws.Cells.Item(14, 17).Value = 1.330788898912
$ws.Cells.Item(14, 18).Value = 11.977100090208
$ws.Cells.Item(14, 19).Value = 0.02693816094070895
$ws.Cells.Item(14, 20).Value = 0.03208321919379618

# Row 15
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 0.790022
$ws.Cells.Item(15, 8).Value = 2.370066
$ws.Cells.Item(15, 9).Value = 0.04193619023665445
$ws.Cells.Item(15, 10).Value = 0.04405545940869034
$ws.Cells.Item(15, 13).Value = 0.9277985
$ws.Cells.Item(15, 14).Value = 1.855597
$ws.Cells.Item(15, 15).Value = 0.3538039290807178
$ws.Cells.Item(15, 16).Value = 0.2674056793968462
$ws.Cells.Item(15, 17).Value = 0.7329812265669999
$ws.Cells.Item(15, 18).Value = 4.397887359402
$ws.Cells.Item(15, 19).Value = 0.01483718887640478
$ws.Cells.Item(15, 20).Value = 0.01178068005432102

# Row 16
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 0.790022
$ws.Cells.Item(16, 8).Value = 2.370066
$ws.Cells.Item(16, 9).Value = 0.04193619023665445
$ws.Cells.Item(16, 10).Value = 0.04405545940869034
$ws.Cells.Item(16, 13).Value = 0.01005766666666667
$ws.Cells.Item(16, 14).Value = 0.030173
$ws.Cells.Item(16, 15).Value = 0.003835360785828855
$ws.Cells.Item(16, 16).Value = 0.004348159414162149
$ws.Cells.Item(16, 17).Value = 0.007945777935333333
$ws.Cells.Item(16, 18).Value = 0.071512001418
$ws.Cells.Item(16, 19).Value = 0.0001608404195407234
$ws.Cells.Item(16, 20).Value = 0.0001915601605731353

